# Cost Analysis workbook update:
#  - Insert a new component row ("Logic Buffer/Scmt. Trig." / "74HC14" / 1 / "3.3 ₺")
#    into the Table1 listing, right after the "Amplifier" row.
#  - Update the running Total Cost shown in the table's totals row from
#    "113.05 ₺" to "116.35 ₺".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a whole new worksheet row at row 5 (pushes the Counter/.../totals rows
# down by one, including everything inside the table).
$ws.Rows.Item(5).Insert()

# Grow the table range by one row so the newly inserted blank row becomes part
# of Table1's data region again (whole-row Insert does not auto-resize it).
$lo.Resize($lo.Range.Resize($lo.Range.Rows.Count + 1))

# Populate the newly-inserted row with the new part.
$ws.Cells.Item(5, 2).Value = "Logic Buffer/Scmt. Trig."
$ws.Cells.Item(5, 3).Value = "74HC14"
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = "3.3 ₺"

# Update the totals-row cost figure (now on row 15 after the insert).
$totalsRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
$ws.Cells.Item($totalsRow, 5).Value = "116.35 ₺"
$lo.ListColumns.Item(4).TotalsRowLabel = "116.35 ₺"
